$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 67: hours, tasks done, additional info
$ws.Range("B67").Value = 8
$ws.Range("C67").Value = "refactoring done, added interaction with tiles by a player "
$ws.Range("D67").Value = "basically an ability to delete or ""break"" tiles from inside of the game"

# Row 68: hours, tasks done
$ws.Range("B68").Value = 5
$ws.Range("C68").Value = "added music, experimented more with replacing variables, textures and sounds "

# Update the view / selection state to match authored file
# (topLeftCell -> B45, selection -> D67)
$excel.ActiveWindow.ScrollRow = 45
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("D67").Select()
